# Update "paises.xlsx" (Pais sheet): refresh COVID case counters for 16
# countries whose totals moved (23 Aug 2020, 18:13 -> 19:30 snapshot), which
# in turn changes the sort-by-total-cases rank/order for a few neighbouring
# countries, and bump the "last updated" timestamp banner.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados ..." banner -------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 23 de Agosto de 2020 a las 19:30"

# --- Rows whose country (column A) moved because of the refreshed ranking -
$nameUpdates = @(
    @{ Row = 52;  Name = "Marruecos" },
    @{ Row = 53;  Name = "Nigeria" },
    @{ Row = 59;  Name = "Etiopia" },
    @{ Row = 60;  Name = "Suiza" },
    @{ Row = 86;  Name = "Libano" },
    @{ Row = 87;  Name = "Sudan" },
    @{ Row = 95;  Name = "Grecia" },
    @{ Row = 96;  Name = "Albania" },
    @{ Row = 97;  Name = "Gabon" },
    @{ Row = 156; Name = "Reunion" },
    @{ Row = 157; Name = "Niger" },
    @{ Row = 202; Name = "Santa Lucia" },
    @{ Row = 203; Name = "Timor Oriental" }
)

foreach ($u in $nameUpdates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.Name
}

# --- Refreshed statistics: Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes (columns B..H) -----
$statUpdates = @(
    @{ Row = 4;   B = 5856346; C = 14918; D = 3152100; E = 2523860; F = 0; G = 212; H = 180386 },
    @{ Row = 13;  B = 397665;  C = 1957;  D = 371179;  E = 15634;   F = 0; G = 60;  H = 10852 },
    @{ Row = 23;  B = 234225;  C = 368;   D = 208950;  E = 15943;   F = 0; G = 1;   H = 9332 },
    @{ Row = 33;  B = 102380;  C = 447;   D = 79501;   E = 22045;   F = 0; G = 15;  H = 834 },
    @{ Row = 52;  B = 52349;   C = 1537;  D = 36343;   E = 15118;   F = 0; G = 30;  H = 888 },
    @{ Row = 53;  B = 51905;   C = 0;     D = 38767;   E = 12141;   F = 0; G = 0;   H = 997 },
    @{ Row = 59;  B = 40671;   C = 1638;  D = 14995;   E = 24998;   F = 0; G = 16;  H = 678 },
    @{ Row = 60;  B = 39903;   C = 276;   D = 34100;   E = 3802;    F = 0; G = 1;   H = 2001 },
    @{ Row = 70;  B = 27969;   C = 61;    D = 23364;   E = 2828;    F = 0; G = 0;   H = 1777 },
    @{ Row = 86;  B = 12698;   C = 507;   D = 3625;    E = 8950;    F = 0; G = 2;   H = 123 },
    @{ Row = 87;  B = 12682;   C = 0;     D = 6492;    E = 5375;    F = 0; G = 0;   H = 815 },
    @{ Row = 95;  B = 8664;    C = 283;   D = 3804;    E = 4618;    F = 0; G = 2;   H = 242 },
    @{ Row = 96;  B = 8427;    C = 152;   D = 4332;    E = 3845;    F = 0; G = 5;   H = 250 },
    @{ Row = 97;  B = 8388;    C = 0;     D = 6734;    E = 1601;    F = 0; G = 0;   H = 53 },
    @{ Row = 156; B = 1209;    C = 92;    D = 692;     E = 511;     F = 0; G = 0;   H = 6 },
    @{ Row = 157; B = 1172;    C = 0;     D = 1083;    E = 20;      F = 0; G = 0;   H = 69 }
)

foreach ($u in $statUpdates) {
    $r = $u.Row
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
    $ws.Cells.Item($r, 6).Value = $u.F
    $ws.Cells.Item($r, 7).Value = $u.G
    $ws.Cells.Item($r, 8).Value = $u.H
}
